$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix typo "avilable" -> "available" in the EEG Recordings paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("avilable", $true, $false, $false, $false, $false, `
    $true, 1, $false, "available", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new sentence before "Specifically, once we recorded..." in the
#    paragraph that starts "In this study, we first used the same strategy."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(". Specifically, once we recorded from a meditator", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ". However, since the number of subjects were less, the strategy was modified to improve matching of demographic data. Specifically, once we recorded from a meditator", `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Append a new sentence at the end of the "Preprocessing:" paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute('within the "preprocessingCodes" folder.', $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    'within the "preprocessingCodes" folder. As mentioned above, 5 subjects are further removed because they had more than 24 bad electrodes, so the final set has 71 subjects.', `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Insert two new paragraphs after the "Preprocessing:" paragraph: a blank
#    paragraph, then a "PowerProject:" paragraph describing the new folder.
#    We clone the (already updated) "Preprocessing:" paragraph twice so the
#    new paragraphs inherit identical paragraph/run formatting, then edit
#    their text in place.
# ---------------------------------------------------------------------------
$prepParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Preprocessing:*") {
        $prepParaIndex = $i
        break
    }
}

$prepPara = $d.Paragraphs.Item($prepParaIndex)
$prepFullRange = $d.Range($prepPara.Range.Start, $prepPara.Range.End)
$prepFormattedText = $prepFullRange.FormattedText

# Clone #1 -> will become the blank spacer paragraph.
$clonePoint1 = $d.Range($prepPara.Range.End, $prepPara.Range.End)
$clonePoint1.FormattedText = $prepFormattedText

# Clone #2 -> will become the "PowerProject:" paragraph.
$clonePoint2 = $d.Range($prepPara.Range.End, $prepPara.Range.End)
$clonePoint2.FormattedText = $prepFormattedText

# --- Clear clone #1 down to an empty paragraph (keep the paragraph mark). ---
$blankPara = $d.Paragraphs.Item($prepParaIndex + 1)
$blankRangeNoMark = $d.Range($blankPara.Range.Start, $blankPara.Range.End - 1)
$blankRangeNoMark.Text = ""

# --- Rewrite clone #2's text as the "PowerProject:" paragraph. -------------
$ppPara = $d.Paragraphs.Item($prepParaIndex + 2)
$ppStart = $ppPara.Range.Start

$headingLen = "Preprocessing: ".Length
$headingRange = $d.Range($ppStart, $ppStart + $headingLen)
$headingRange.Text = "PowerProject: "

$ppParaAfterHeading = $d.Paragraphs.Item($prepParaIndex + 2)
$newHeadingEnd = $ppStart + "PowerProject: ".Length
$bodyRange = $d.Range($newHeadingEnd, $ppParaAfterHeading.Range.End - 1)
$bodyRange.Text = "This folder contains programs to display the power for different protocols. More details can be found in the ReadMe file within that folder."
